# Update header labels on every worksheet (April, May, June):
#   C1: "Consumption Days" -> "Consumption Period"
#   D1: "Usage (%)"        -> "Utilisation (%)"

$wb = $excel.ActiveWorkbook

foreach ($ws in $wb.Worksheets) {
    $ws.Range("C1").Value = "Consumption Period"
    $ws.Range("D1").Value = "Utilisation (%)"
}
